# Update quarterly financial figures for CHL (Income Statement, Balance Sheet, Cash Flow)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHL")

# Income Statement
$ws.Range("D8").Value = 58151800
$ws.Range("E8").Value = 52187300
$ws.Range("F8").Value = 57712300
$ws.Range("G8").Value = 50173000
$ws.Range("H8").Value = 54963800
$ws.Range("I8").Value = 47883600
$ws.Range("J8").Value = 51304000
$ws.Range("D9").Value = 10566300
$ws.Range("E9").Value = 9681700
$ws.Range("F9").Value = 11357800
$ws.Range("G9").Value = 10251100
$ws.Range("H9").Value = 11745300
$ws.Range("I9").Value = 10230600
$ws.Range("J9").Value = 9305000
$ws.Range("D10").Value = 47585400
$ws.Range("E10").Value = 42505700
$ws.Range("F10").Value = 46354500
$ws.Range("G10").Value = 39921800
$ws.Range("H10").Value = 43218500
$ws.Range("I10").Value = 37653000
$ws.Range("J10").Value = 41999000
$ws.Range("D12").Value = 3459600
$ws.Range("E12").Value = 5642000
$ws.Range("E14").Value = 1868900
$ws.Range("G14").Value = 1070900
$ws.Range("D15").Value = 11309000
$ws.Range("E15").Value = 11510500
$ws.Range("F15").Value = 10794700
$ws.Range("G15").Value = 10443600
$ws.Range("H15").Value = 10124400
$ws.Range("I15").Value = 10264200
$ws.Range("J15").Value = 10083700
$ws.Range("D17").Value = 47809800
$ws.Range("E17").Value = 44447500
$ws.Range("F17").Value = 47624300
$ws.Range("G17").Value = 42462000
$ws.Range("H17").Value = 45149300
$ws.Range("I17").Value = 41909800
$ws.Range("J17").Value = 42003100
$ws.Range("D18").Value = 10342000
$ws.Range("E18").Value = 7739900
$ws.Range("F18").Value = 10088000
$ws.Range("G18").Value = 7710900
$ws.Range("H18").Value = 9814500
$ws.Range("I18").Value = 5973800
$ws.Range("J18").Value = 9300900
$ws.Range("D20").Value = 2157000
$ws.Range("E20").Value = 2165300
$ws.Range("F20").Value = 1991800
$ws.Range("G20").Value = 2033400
$ws.Range("H20").Value = 1880800
$ws.Range("I20").Value = 4326900
$ws.Range("J20").Value = 1730000
$ws.Range("D21").Value = 23807900
$ws.Range("E21").Value = 10621000
$ws.Range("F21").Value = 22874600
$ws.Range("G21").Value = 10063500
$ws.Range("H21").Value = 21819700
$ws.Range("I21").Value = 10979100
$ws.Range("J21").Value = 21114600
$ws.Range("D23").Value = 12498900
$ws.Range("E23").Value = 9905200
$ws.Range("F23").Value = 12079800
$ws.Range("G23").Value = 9744300
$ws.Range("H23").Value = 11695300
$ws.Range("I23").Value = 10300700
$ws.Range("J23").Value = 11030900
$ws.Range("D24").Value = 2742200
$ws.Range("E24").Value = 2234300
$ws.Range("F24").Value = 2770500
$ws.Range("G24").Value = 2587800
$ws.Range("H24").Value = 2699000
$ws.Range("I24").Value = 2694100
$ws.Range("J24").Value = 2512000
$ws.Range("D26").Value = 9756800
$ws.Range("E26").Value = 7670900
$ws.Range("F26").Value = 9309300
$ws.Range("G26").Value = 7156500
$ws.Range("H26").Value = 8996300
$ws.Range("I26").Value = 7606600
$ws.Range("J26").Value = 8518900
$ws.Range("D27").Value = 9741800
$ws.Range("E27").Value = 7658500
$ws.Range("F27").Value = 9301600
$ws.Range("G27").Value = 7148800
$ws.Range("H27").Value = 8989500
$ws.Range("I27").Value = 7599000
$ws.Range("J27").Value = 8509200
$ws.Range("D32").Value = -2157000
$ws.Range("E32").Value = -2165300
$ws.Range("F32").Value = -1991800
$ws.Range("G32").Value = -2033400
$ws.Range("H32").Value = -1880800
$ws.Range("I32").Value = -4326900
$ws.Range("J32").Value = -1730000
$ws.Range("D33").Value = 9741800
$ws.Range("E33").Value = 7658500
$ws.Range("F33").Value = 9301600
$ws.Range("G33").Value = 7148800
$ws.Range("H33").Value = 8989500
$ws.Range("I33").Value = 7599000
$ws.Range("J33").Value = 8509200
$ws.Range("D35").Value = 9741800
$ws.Range("E35").Value = 7658500
$ws.Range("F35").Value = 9301600
$ws.Range("G35").Value = 7148800
$ws.Range("H35").Value = 8989500
$ws.Range("I35").Value = 7599000
$ws.Range("J35").Value = 8509200
# Balance Sheet
$ws.Range("D41").Value = 12281100
$ws.Range("E41").Value = 17903600
$ws.Range("F41").Value = 12319200
$ws.Range("G41").Value = 13418200
$ws.Range("H41").Value = 14217700
$ws.Range("I41").Value = 11849400
$ws.Range("J41").Value = 14964600
$ws.Range("D42").Value = 57289100
$ws.Range("E42").Value = 51201600
$ws.Range("F42").Value = 58631000
$ws.Range("G42").Value = 54495300
$ws.Range("H42").Value = 54884800
$ws.Range("I42").Value = 50830000
$ws.Range("J42").Value = 52932800
$ws.Range("D43").Value = 10898800
$ws.Range("E43").Value = 8473300
$ws.Range("F43").Value = 16702400
$ws.Range("G43").Value = 15317100
$ws.Range("H43").Value = 6075500
$ws.Range("I43").Value = 6666900
$ws.Range("J43").Value = 5091800
$ws.Range("D44").Value = 1579200
$ws.Range("E44").Value = 1517000
$ws.Range("F44").Value = 1084400
$ws.Range("G44").Value = 1310800
$ws.Range("H44").Value = 1031000
$ws.Range("I44").Value = 1483200
$ws.Range("J44").Value = 988000
$ws.Range("D45").Value = 3230100
$ws.Range("E45").Value = 3746300
$ws.Range("F45").Value = 2662600
$ws.Range("G45").Value = 2522700
$ws.Range("H45").Value = 1747800
$ws.Range("I45").Value = 1698100
$ws.Range("J45").Value = 2238300
$ws.Range("D46").Value = 85278300
$ws.Range("E46").Value = 82841900
$ws.Range("F46").Value = 91399600
$ws.Range("G46").Value = 87064000
$ws.Range("H46").Value = 77956800
$ws.Range("I46").Value = 72527500
$ws.Range("J46").Value = 76215500
$ws.Range("D47").Value = 20299800
$ws.Range("E47").Value = 19670700
$ws.Range("F47").Value = 18975000
$ws.Range("G47").Value = 18413800
$ws.Range("H47").Value = 26342200
$ws.Range("I47").Value = 25626400
$ws.Range("J47").Value = 10617700
$ws.Range("D48").Value = 108227400
$ws.Range("E48").Value = 107766600
$ws.Range("F48").Value = 107442900
$ws.Range("G48").Value = 105698900
$ws.Range("H48").Value = 101935300
$ws.Range("I48").Value = 99975400
$ws.Range("J48").Value = 98354300
$ws.Range("D49").Value = 5587600
$ws.Range("E49").Value = 5500700
$ws.Range("F49").Value = 5485800
$ws.Range("G49").Value = 5498700
$ws.Range("H49").Value = 5359200
$ws.Range("I49").Value = 5359200
$ws.Range("J49").Value = 5348100
$ws.Range("D52").Value = 12886300
$ws.Range("E52").Value = 10117000
$ws.Range("F52").Value = 10427000
$ws.Range("G52").Value = 9055200
$ws.Range("H52").Value = 9660200
$ws.Range("I52").Value = 8425400
$ws.Range("J52").Value = 8750600
$ws.Range("D54").Value = 232279500
$ws.Range("E54").Value = 225896800
$ws.Range("F54").Value = 233730300
$ws.Range("G54").Value = 225730700
$ws.Range("H54").Value = 221253600
$ws.Range("I54").Value = 211913900
$ws.Range("J54").Value = 199286100
$ws.Range("D57").Value = 31683200
$ws.Range("E57").Value = 35094800
$ws.Range("F57").Value = 35934700
$ws.Range("G57").Value = 37405900
$ws.Range("H57").Value = 34858700
$ws.Range("I57").Value = 36245300
$ws.Range("J57").Value = 30427900
$ws.Range("F58").Value = 742100
$ws.Range("G58").Value = 741800
$ws.Range("J58").Value = 10100
$ws.Range("D59").Value = 46410800
$ws.Range("E59").Value = 43559800
$ws.Range("F59").Value = 45055000
$ws.Range("G59").Value = 41457900
$ws.Range("H59").Value = 42959400
$ws.Range("I59").Value = 38113800
$ws.Range("J59").Value = 35153400
$ws.Range("D60").Value = 78093900
$ws.Range("E60").Value = 78654600
$ws.Range("F60").Value = 81731800
$ws.Range("G60").Value = 79605500
$ws.Range("H60").Value = 77818000
$ws.Range("I60").Value = 74359000
$ws.Range("J60").Value = 65591400
$ws.Range("H61").Value = 741500
$ws.Range("I61").Value = 741300
$ws.Range("J61").Value = 741000
$ws.Range("D62").Value = 577200
$ws.Range("E62").Value = 482300
$ws.Range("F62").Value = 393600
$ws.Range("G62").Value = 366100
$ws.Range("H62").Value = 219600
$ws.Range("I62").Value = 221700
$ws.Range("J62").Value = 149600
$ws.Range("D66").Value = 79167700
$ws.Range("E66").Value = 79618600
$ws.Range("F66").Value = 82595400
$ws.Range("G66").Value = 80434200
$ws.Range("H66").Value = 79236000
$ws.Range("I66").Value = 75772100
$ws.Range("J66").Value = 66795300
$ws.Range("D72").Value = 93407000
$ws.Range("E72").Value = 86616800
$ws.Range("F72").Value = 91427500
$ws.Range("G72").Value = 85526000
$ws.Range("H72").Value = 82326600
$ws.Range("I72").Value = 76486200
$ws.Range("J72").Value = 72928500
$ws.Range("D76").Value = 153111800
$ws.Range("E76").Value = 146278200
$ws.Range("F76").Value = 151135000
$ws.Range("G76").Value = 145296500
$ws.Range("H76").Value = 142017700
$ws.Range("I76").Value = 136141800
$ws.Range("J76").Value = 132490800
# Cash Flow Statement
$ws.Range("D81").Value = 9741800
$ws.Range("E81").Value = 7658500
$ws.Range("F81").Value = 9301600
$ws.Range("G81").Value = 7148800
$ws.Range("H81").Value = 8989500
$ws.Range("I81").Value = 7599000
$ws.Range("J81").Value = 8509200
$ws.Range("D89").Value = 17559400
$ws.Range("E89").Value = 15916800
$ws.Range("F89").Value = 20519900
$ws.Range("G89").Value = 16295700
$ws.Range("H89").Value = 21356100
$ws.Range("I89").Value = 15710800
$ws.Range("J89").Value = 19178700
$ws.Range("D91").Value = -11798600
$ws.Range("E91").Value = -15978900
$ws.Range("F91").Value = -12666500
$ws.Range("G91").Value = -15647600
$ws.Range("H91").Value = -12284500
$ws.Range("J91").Value = -10353800
$ws.Range("D94").Value = -19085800
$ws.Range("E94").Value = 1663700
$ws.Range("F94").Value = -17474200
$ws.Range("G94").Value = -14001100
$ws.Range("H94").Value = -14868000
$ws.Range("I94").Value = -9892900
$ws.Range("J94").Value = -11291600
$ws.Range("D100").Value = -4136500
$ws.Range("E100").Value = -11962600
$ws.Range("F100").Value = -4100000
$ws.Range("G100").Value = -3130300
$ws.Range("H100").Value = -4135600
$ws.Range("I100").Value = -9657600
$ws.Range("J100").Value = -3181300
$ws.Range("D101").Value = 40400
$ws.Range("E101").Value = -33500
$ws.Range("F101").Value = -44700
$ws.Range("G101").Value = 36200
$ws.Range("H101").Value = 15900
$ws.Range("I101").Value = 36700
$ws.Range("J101").Value = -7900
$ws.Range("D102").Value = -5622500
$ws.Range("E102").Value = 5584400
$ws.Range("F102").Value = -1099000
$ws.Range("G102").Value = -799500
$ws.Range("H102").Value = 2368300
$ws.Range("I102").Value = -3803000
$ws.Range("J102").Value = 4697900
